$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.881.14"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.034.09"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.93"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.47"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0817"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "2.335.11"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.51"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.31"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.761"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "2.035.75"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "37.823.61"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -6.61%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.09"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.99"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.86"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.26"
$ws.Range("E32").Value = "  +9.18%  "
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0605"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.33"
$ws.Range("E36").Value = "  +4.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.28"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.60"
$ws.Range("E40").Value = "  +6.36%  "
$ws.Range("D41").Value = "1.530.83"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.23"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0914"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.98"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.96"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "2.225.36"
$ws.Range("E51").Value = "  -0.79%  "
